$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.117.30'
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("D3").Value = '2.471.46'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  +2.66%  '
$ws.Range("D9").Value = '2.478.30'
$ws.Range("E9").Value = '  +1.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0966'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.331'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.51%  '
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("D14").Value = '2.905.54'
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("D15").Value = '56.168.12'
$ws.Range("E15").Value = '  +2.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000136'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.44%  '
$ws.Range("D18").Value = '2.477.65'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '317.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +7.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.412'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.52%  '
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("E27").Value = '  +3.64%  '
$ws.Range("D28").Value = '2.583.19'
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.43%  '
$ws.Range("D30").Value = '0.0₃0787'
$ws.Range("E30").Value = '  +8.99%  '
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("E34").Value = '  +4.46%  '
$ws.Range("E35").Value = '  +3.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.860'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.54%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0555'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.53%  '
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("E44").Value = '  +6.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.10%  '
$ws.Range("E46").Value = '  +4.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '257.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.81%  '
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("E49").Value = '  +4.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.02%  '
$ws.Range("D51").Value = '1.869.23'
$ws.Range("E51").Value = '  -3.95%  '
